$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 6) with data for the "RijndaelEncryption" target application,
# used in validating SharpChecker.
$ws.Range("A6").Value = "RijndaelEncryption"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 100
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = "Encrypted"
$ws.Range("G6").Value = 2

# Update the selection to reflect the new active cell after data entry.
$ws.Range("C7").Select()
